# Mark the remaining "WeatherReport" checklist rows (11-16, i.e. Q.9-Q.14)
# as Done in columns B:E (Pandas/Excel/PowerBi/Sql), matching the other
# already-completed rows above them. Then leave the on-screen selection on
# A23 for both sheets (where the user had scrolled to).

$wb = $excel.ActiveWorkbook

$wsReport = $wb.Worksheets.Item("WeatherReport")
$wsExplain = $wb.Worksheets.Item("WeatherReportExplanations")

# Rows 11..16 correspond to questions 9..14 that were still blank.
$doneRange = $wsReport.Range("B11:E16")
$doneRange.Value = "Done"

# Restore the saved cursor/selection positions recorded in the file.
$wsReport.Activate()
$wsReport.Range("A23").Select()

$wsExplain.Activate()
$wsExplain.Application.ActiveWindow.ScrollRow = 23
$wsExplain.Range("A23").Select()
